$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace "J.Plaggenberg" with "J.PLAGGENBERG" in columns B and E (s2cDNAPreparer, libraryPreparer)
# for all data rows (2 through 42).
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 2).Value = "J.PLAGGENBERG"
    $ws.Cells.Item($r, 5).Value = "J.PLAGGENBERG"
}

# Update the active selection to M3 (as recorded in the saved view state).
$ws.Range("M3").Select() | Out-Null
